$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet ---
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data rows
$wsForecast.Range("A2").Value = 45319.99999999999
$wsForecast.Range("B2").Value = 108
$wsForecast.Range("C2").Value = -107.1758655516373
$wsForecast.Range("D2").Value = 316.3282464665403
$wsForecast.Range("A3").Value = 45326.99999999999
$wsForecast.Range("B3").Value = 111
$wsForecast.Range("C3").Value = -95.99015196281276
$wsForecast.Range("D3").Value = 332.9824563036667
$wsForecast.Range("A4").Value = 45333.99999999999
$wsForecast.Range("B4").Value = 114
$wsForecast.Range("C4").Value = -97.60664833951438
$wsForecast.Range("D4").Value = 310.2039932699041
$wsForecast.Range("A5").Value = 45347.99999999999
$wsForecast.Range("B5").Value = 120
$wsForecast.Range("C5").Value = -95.17215383496055
$wsForecast.Range("D5").Value = 342.2846932221948
$wsForecast.Range("A6").Value = 45361.99999999999
$wsForecast.Range("B6").Value = 126
$wsForecast.Range("C6").Value = -86.40575005915483
$wsForecast.Range("D6").Value = 345.3815117004606
$wsForecast.Range("A7").Value = 45368.99999999999
$wsForecast.Range("B7").Value = 129
$wsForecast.Range("C7").Value = -87.77328677876986
$wsForecast.Range("D7").Value = 326.4300043824768
$wsForecast.Range("A8").Value = 45375.99999999999
$wsForecast.Range("B8").Value = 133
$wsForecast.Range("C8").Value = -90.48356752404163
$wsForecast.Range("D8").Value = 354.5412984296087
$wsForecast.Range("A9").Value = 45382.99999999999
$wsForecast.Range("B9").Value = 136
$wsForecast.Range("C9").Value = -67.79417614471828
$wsForecast.Range("D9").Value = 352.491038192185
$wsForecast.Range("A10").Value = 45396.99999999999
$wsForecast.Range("B10").Value = 142
$wsForecast.Range("C10").Value = -70.87530017881009
$wsForecast.Range("D10").Value = 351.7587960975926
$wsForecast.Range("A11").Value = 45403.99999999999
$wsForecast.Range("B11").Value = 145
$wsForecast.Range("C11").Value = -71.05343410917703
$wsForecast.Range("D11").Value = 355.4948326213589
$wsForecast.Range("A12").Value = 45424.99999999999
$wsForecast.Range("B12").Value = 154
$wsForecast.Range("C12").Value = -70.39813088607815
$wsForecast.Range("D12").Value = 380.9920553893932
$wsForecast.Range("A13").Value = 45438.99999999999
$wsForecast.Range("B13").Value = 160
$wsForecast.Range("C13").Value = -59.96693612915158
$wsForecast.Range("D13").Value = 391.4548165901609
$wsForecast.Range("A14").Value = 45445.99999999999
$wsForecast.Range("B14").Value = 163
$wsForecast.Range("C14").Value = -50.01750283161895
$wsForecast.Range("D14").Value = 391.5124691240165
$wsForecast.Range("A15").Value = 45452.99999999999
$wsForecast.Range("B15").Value = 166
$wsForecast.Range("C15").Value = -48.22247976157409
$wsForecast.Range("D15").Value = 379.928894099297
$wsForecast.Range("A16").Value = 45459.99999999999
$wsForecast.Range("B16").Value = 169
$wsForecast.Range("C16").Value = -45.32464945017347
$wsForecast.Range("D16").Value = 395.5344669996625
$wsForecast.Range("A17").Value = 45487.99999999999
$wsForecast.Range("B17").Value = 182
$wsForecast.Range("C17").Value = -25.00031199036649
$wsForecast.Range("D17").Value = 384.3277054466652
$wsForecast.Range("A18").Value = 45494.99999999999
$wsForecast.Range("B18").Value = 185
$wsForecast.Range("C18").Value = -29.21981632801959
$wsForecast.Range("D18").Value = 395.1264815249154
$wsForecast.Range("A19").Value = 45501.99999999999
$wsForecast.Range("B19").Value = 188
$wsForecast.Range("C19").Value = -25.33074361294591
$wsForecast.Range("D19").Value = 405.4791087918284
$wsForecast.Range("A20").Value = 45515.99999999999
$wsForecast.Range("B20").Value = 194
$wsForecast.Range("C20").Value = -12.7000942181184
$wsForecast.Range("D20").Value = 405.7465605881011
$wsForecast.Range("A21").Value = 45522.99999999999
$wsForecast.Range("B21").Value = 197
$wsForecast.Range("C21").Value = -23.95734099712536
$wsForecast.Range("D21").Value = 415.6080942633634
$wsForecast.Range("A22").Value = 45529.99999999999
$wsForecast.Range("B22").Value = 200
$wsForecast.Range("C22").Value = -5.200370576845031
$wsForecast.Range("D22").Value = 418.2411818508918
$wsForecast.Range("A23").Value = 45536.99999999999
$wsForecast.Range("B23").Value = 203
$wsForecast.Range("C23").Value = -5.040012542975731
$wsForecast.Range("D23").Value = 414.6174623746853
$wsForecast.Range("A24").Value = 45543.99999999999
$wsForecast.Range("B24").Value = 206
$wsForecast.Range("C24").Value = -6.796863510577356
$wsForecast.Range("D24").Value = 418.1823586524091
$wsForecast.Range("A25").Value = 45550.99999999999
$wsForecast.Range("B25").Value = 209
$wsForecast.Range("C25").Value = -12.00135328816653
$wsForecast.Range("D25").Value = 413.0993373758053
$wsForecast.Range("A26").Value = 45557.99999999999
$wsForecast.Range("B26").Value = 212
$wsForecast.Range("C26").Value = -5.242623053972896
$wsForecast.Range("D26").Value = 430.0541637826004
$wsForecast.Range("A27").Value = 45564.99999999999
$wsForecast.Range("B27").Value = 215
$wsForecast.Range("C27").Value = -4.468009335481562
$wsForecast.Range("D27").Value = 434.1022269227804
$wsForecast.Range("A28").Value = 45571.99999999999
$wsForecast.Range("B28").Value = 219
$wsForecast.Range("C28").Value = 2.178156858745199
$wsForecast.Range("D28").Value = 445.2260229859896
$wsForecast.Range("A29").Value = 45578.99999999999
$wsForecast.Range("B29").Value = 222
$wsForecast.Range("C29").Value = 7.405013391596692
$wsForecast.Range("D29").Value = 453.0920283038919
$wsForecast.Range("A30").Value = 45585.99999999999
$wsForecast.Range("B30").Value = 225
$wsForecast.Range("C30").Value = 2.998071937823639
$wsForecast.Range("D30").Value = 444.5467926378341
$wsForecast.Range("A31").Value = 45592.99999999999
$wsForecast.Range("B31").Value = 228
$wsForecast.Range("C31").Value = 13.18398827356134
$wsForecast.Range("D31").Value = 460.2993687402687
$wsForecast.Range("A32").Value = 45599.99999999999
$wsForecast.Range("B32").Value = 231
$wsForecast.Range("C32").Value = 6.68090137041275
$wsForecast.Range("D32").Value = 435.4098220582038

# Apply date/time number format to column A (rows 2-32), matching the other sheets
$wsForecast.Range("A2:A32").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Move the new sheet to the end (after "Monthly Trend") ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsForecast.Move([Type]::Missing, $wsMonthly)

